$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.508.71"
$ws.Range("E2").Value = "  -1.01%  "
$ws.Range("D3").Value = "1.914.33"
$ws.Range("E3").Value = "  -1.45%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.20%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4787"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.13%  "
$ws.Range("E8").Value = "  -3.99%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06701"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.66%  "
$ws.Range("E10").Value = "  -3.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "101.79"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07713"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.24%  "
$ws.Range("D13").Value = "1.917.93"
$ws.Range("E13").Value = "  -1.28%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.200"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.72%  "
$ws.Range("E15").Value = "  -3.88%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "262.86"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.22%  "
$ws.Range("D17").Value = "30.504.91"
$ws.Range("E17").Value = "  -1.02%  "
$ws.Range("E18").Value = "  +0.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007458"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.23%  "
$ws.Range("E20").Value = "  -3.61%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.399"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.003"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.292"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.43%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.350"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.82%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "166.89"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "19.13"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.70%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.062"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.52%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.387"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.32%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.677"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.56%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09977"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.513"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.65%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.236"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.00%  "
$ws.Range("E33").Value = "  -2.51%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7251"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.58%  "
$ws.Range("E35").Value = "  -4.45%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.723"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01914"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.617"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.51%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.241"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "74.66"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.67%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.965"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.42%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8627"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.72%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "105.46"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.42%  "
$ws.Range("E44").Value = "  -3.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.002"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.37%  "
$ws.Range("E46").Value = "  -4.42%  "
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "928.34"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.50%  "
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1203"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.41%  "
$ws.Range("E49").Value = "  -3.59%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05754"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.44%  "
